$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
$c1 = $cs.Colors(1)
Write-Host "before RGB:" $c1.RGB
$c1.RGB = 255
Write-Host "after RGB:" $c1.RGB
